# Junction_Flooding_34: round J1..J33 readings in row 5 to 2 decimal places
# ("custom accuracy") and drop the extra data row 6 (reducing the sample
# down, part of the "데이터 1000개" trim across the workbook set).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Round row 5 (B5:AH5) values to 2 decimal places -------------------
$lastCol = 34   # column AH
for ($c = 2; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(5, $c)
    $cell.Value = [Math]::Round([double]$cell.Value2, 2)
}

# --- 2. Delete row 6 entirely ----------------------------------------------
$ws.Rows(6).Delete()

# --- 3. Re-fit column widths B:AH to the now-shorter numbers --------------
# ColumnWidth is in characters; the stored <col width> ends up
# ColumnWidth + 5/6 (the default-font padding), so back that constant off
# to land on the exact integer widths Excel's own AutoFit produced.
$colWidths = @{
    2 = 7; 3 = 7; 4 = 7; 5 = 8; 6 = 8; 7 = 7; 8 = 8; 9 = 8; 10 = 7;
    11 = 7; 12 = 7; 13 = 7; 14 = 7; 15 = 7; 16 = 7; 17 = 7; 18 = 7; 19 = 7;
    20 = 8; 21 = 8; 22 = 7; 23 = 8; 24 = 7; 25 = 7; 26 = 8; 27 = 7; 28 = 7;
    29 = 7; 30 = 7; 31 = 7; 32 = 8; 33 = 7; 34 = 7
}
foreach ($c in $colWidths.Keys) {
    $ws.Columns($c).ColumnWidth = $colWidths[$c] - (5/6)
}
